$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 23 (shift down, copying format from the row above),
# pushing the existing rows 23-27 down to 24-28.
$ws.Rows.Item(23).Insert(-4121, 0)

# Populate the new row with the MPLS data point: date issued 7/29/2020
# (serial 44041), no EO #, purpose describing the indoor dining closure.
$ws.Range("A23").Value = 44041
$ws.Range("C23").Value = "Mayor Frey orders indoor dining closed 8/1"

# Leave the active selection on C23, matching the post-edit UI state.
$ws.Range("C23").Select()
